# Updated cryptos list — refresh Price / Volume(1h) figures, and swap the
# InjectiveProtocol / Monero rows (46/47) to their new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 46 & 47 swap places (Monero now ranks above InjectiveProtocol) ---
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'124.94"
$ws.Range("E46").Value = "  +3.08%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'25.52"
$ws.Range("E47").Value = "  +1.34%  "

# --- Refreshed Price (D) / Volume(1h) (E) values for the remaining rows ---
$ws.Range("D2").Value = "63.533.82"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.083.33"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'544.26"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "'140.01"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.078.00"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "'0.505"
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'6.36"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("D13").Value = "'34.99"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").Value = "3.584.68"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "63.558.84"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'0.113"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "3.082.32"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").Value = "'475.32"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("D22").Value = "'0.701"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").Value = "'78.88"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'12.28"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").Value = "'7.97"
$ws.Range("E28").Value = "  -6.46%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'26.23"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("E31").Value = "  -3.86%  "
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("D33").Value = "'58.16"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("E34").Value = "  -7.69%  "
$ws.Range("D35").Value = "'5.41"
$ws.Range("E35").Value = "  +4.98%  "
$ws.Range("D36").Value = "'491.77"
$ws.Range("E36").Value = "  -5.30%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "3.247.32"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").Value = "'0.0405"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("D42").Value = "'8.13"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").Value = "'0.254"
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("D48").Value = "'2.04"
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("D49").Value = "0.0₃0529"
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("E51").Value = "  +0.18%  "
